# Reading data from excel
# Module: AppContent
# Providing inputs to upload file through data provider by reading data
# from excel sheet.
#
# The "AppContent" sheet used to be a 4-column table (TestData header +
# CaseFlag/Expected result/Result columns). It becomes a single-column
# list of upload-file test data: the existing "upload_bmp.exe" row stays,
# and two new rows are appended for "upload_gif.exe" and
# "upload_jpeg.exe". The now-unused helper columns B:D are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AppContent")
$ws.Activate()

# New test data rows (A3 was blank before, A4 is brand new).
$ws.Range("A3").Value = "./src/com/autoitfiles/upload_gif.exe"
$ws.Range("A4").Value = "./src/com/autoitfiles/upload_jpeg.exe"

# A4 is a brand-new row - give it the same look (border/font/alignment)
# as the other data rows by copying A3's formatting onto it.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows("4:4").RowHeight = 15.75

# The sheet is now single-column data (TestData), so drop the old
# CaseFlag / Expected result / Result columns.
$ws.Columns("B:D").Delete()

$ws.Range("A4").Select()
